# Swap the data values between row 3 and row 4 for the columns that differ
# (A, B, E, F, G, H, Q, R, Z, AB) while leaving the other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $columns) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"

    $val3 = $ws.Range($addr3).Value2
    $val4 = $ws.Range($addr4).Value2

    $ws.Range($addr3).Value = $val4
    $ws.Range($addr4).Value = $val3
}
